# La Penita workbook update: column E ("04dec2025") gains real counts for
# several empadronadores, on both the numeric "crosstab" sheet and its
# text-mirrored "annot" sheet.
#
# row -> new value for column E (5)
#   4  CHERO JUAREZ ANYELA TATIANA      0  -> 8
#   5  GARAVITO LEON IVONNE LISSETH     0  -> 8
#   8  NIÑO GUERRERO ANYELA MELINA      0  -> 10
#   9  PANTA MONZON SHIRLEY MARIBEL     0  -> 15
#  10  PEREZ VEGA ANA YSABEL           30  -> 32
#  11  TIMOTEO BAYONA SHARYN LISSETH    0  -> 8
#  12  TIZON NUÑEZ FRESIA YAMILI        0  -> 10
#  13  VALLE SILVA SUTMMER ORFELINDA    0  -> 14
#  14  ZAPATA ZETA ROSA ARACELI         0  -> 13

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("crosstab")
$ws2 = $wb.Worksheets.Item("annot")

$updates = [ordered]@{
    4  = 8
    5  = 8
    8  = 10
    9  = 15
    10 = 32
    11 = 8
    12 = 10
    13 = 14
    14 = 13
}

# "annot" stores every figure as text (it mirrors "crosstab" for display).
# Assigning a numeric-looking string straight to .Value would make Excel
# re-interpret it as a Number, so instead we stage the text in a scratch
# cell well outside the used range (A1:K14), force it to Text with a "@"
# number format, and PasteSpecial just the values into the target cell.
# That carries over the Text type without leaving any NumberFormat/style
# change behind on the destination cells.
$scratch = $ws2.Range("Z1")
$scratch.NumberFormat = "@"

foreach ($row in $updates.Keys) {
    $val = $updates[$row]

    # Numeric crosstab sheet: plain numeric assignment.
    $ws1.Cells.Item($row, 5).Value = $val

    # Text-mirrored annot sheet: paste-special values only.
    $scratch.Value = [string]$val
    $scratch.Copy()
    $ws2.Cells.Item($row, 5).PasteSpecial(-4163)
}

$scratch.Clear()
